# Auto-generated PowerShell COM-interop script
# Applies the LOQ4240.docx content reshuffle described by the commit diff.
$d = $word.ActiveDocument

function Set-ParaText([int]$index, [string]$newText) {
    $p = $d.Paragraphs.Item($index)
    $rng = $p.Range
    $rng.End = $rng.End - 1   # exclude the paragraph mark
    $rng.Text = $newText
}

$BR = [char]11   # manual line break (renders as <w:br/>)

# --- Paragraphs 6, 7, 9, 12, 14, 19: simple single-run text swaps (style unchanged) ---
Set-ParaText 6 "Conceitos fundamentais de administração e noções básicas de marketing e Gestão de Pessoas."
Set-ParaText 7 "Fundamental concepts of administration and basic notions of marketing and People Management."
Set-ParaText 9 "Introduzir os conceitos fundamentais de administração, de configurações de uma organização, de marketing e Gestão de Pessoas, de forma genérica. A disciplina privilegia a discussão dos fundamentos das diversas abordagens e linhas de pensamento administrativo, sob a ótica da engenharia."
Set-ParaText 12 "To Introduce the fundamental concepts of administration, configurations of an organization, marketing and People Management, in a generic way. The subject privileges the discussion of the fundamentals of the different approaches and lines of administrative thought, from the perspective of engineering."
Set-ParaText 14 "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."
Set-ParaText 19 "11079086 - Herlandí de Souza Andrade"

# --- Paragraph 11: Programa (PT) list, now living in the 'Programa resumido' slot ---
$p11text = "1. Elementos de organizações de alto desempenho: aprendizagem organizacional, modelo da competência e capacitações dinâmicas." + $BR + "2. Introdução à Gestão de Pessoas" + $BR + "3. Noções básicas de Marketing" + $BR + "4. Desenvolvimento de atividade prática extensionista junto à micro e pequenos empreendedores da região (componente curricular: plano de marketing)" + $BR + "5. Visita (viagem didática complementar) a uma empresa para conhecer e entender os diferentes processos organizacionais."
Set-ParaText 11 $p11text

# --- Paragraph 17 (Avaliacao bullet list): rotate label values + absorb bibliography ---
# Método's value becomes the old Critério value;
# Critério's value becomes the old Norma de recuperação value;
# Norma de recuperação's value becomes the full bibliography block (ex-paragraph 19).
$p17 = $d.Paragraphs.Item(17)
$p17rng = $p17.Range

# Replace 'Método:' value (scoped to paragraph 17 to avoid cross-paragraph matches)
$p17rng.Find.Execute("Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras.", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_METODO_TOKEN@@", 2) | Out-Null
$p17rng = $p17.Range
$p17rng.Find.Execute("Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_CRITERIO_TOKEN@@", 2) | Out-Null
$p17rng = $p17.Range
$p17rng.Find.Execute("NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_NORMA_TOKEN@@", 2) | Out-Null

# Second pass: placeholders -> final values (placeholders are unique, so no collisions)
$p17rng = $p17.Range
$p17rng.Find.Execute("@@PLACEHOLDER_METODO_TOKEN@@", $true, $false, $false, $false, $false, $true, 1, $false, "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas", 2) | Out-Null
$p17rng = $p17.Range
$p17rng.Find.Execute("@@PLACEHOLDER_CRITERIO_TOKEN@@", $true, $false, $false, $false, $false, $true, 1, $false, "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação", 2) | Out-Null

# Norma de recuperação's value -> entire bibliography block (built with explicit breaks)
$bibText = "Chiavenato, I. Gestão de Pessoas. 4 ed. São Paulo: Manole, 2014." + $BR + $BR + "Chiavenato, I. Recursos Humanos: o capital humano das organizações. 10 ed. Rio de Janeiro, Campus, 2015." + $BR + $BR + "ROBBINS, S. P.; DECENZO, D. A.; WOLTER, R. Fundamentos de Gestão de Pessoas. São Paulo, saraiva, 2013." + $BR + $BR + "KOTLER, P. - ARMSTRONG, G. Princípios De Marketing. 15 ed. São Paulo: Pearson, 2014." + $BR + $BR + "KOTLER, P.; KELLER, K. L. Administração De Marketing. 15 ed. São Paulo: Pearson, 2019." + $BR + $BR + "CHIAVENATO, I. Introdução À Teoria Geral da Administração. 9 ed. São Paulo: Manole, 2014. " + $BR + $BR + "MAXIMIANO, A. C. A. Teoria Geral da Administração: da Revolução Urbana À Revolução Digital. 8 ed. São Paulo: Atlas, 2017." + $BR + $BR + "GUERRINI, F. M.; ESCRIÇÃO FILHO, E.; ROSIM, D. Administração Para Engenheiros. Rio de Janeiro: Campus, 2016." + $BR + $BR + "CHIAVENATO, I. Administração Para Não Administradores: a Gestão de Negócios Ao Alcance de Todos. 2 ed. São Paulo: Manole, 2011." + $BR + $BR + "SILVA, M. M. L. Administração para Estudantes e Profissionais de Áreas Técnicas. São Paulo: Brasport, 2018." + $BR + $BR + "BOLMAN, L.G.; DEAL, T.E. Reframing organizations. San Francisco, John Wiley, 2013" + $BR + $BR + "KOTLER, P.. O Marketing sem segredos. 1 ed. Porto Alegre. Bookman, 2005" + $BR + $BR + "MINTZBERG, H. Criando organizações eficazes. 2 ed. São Paulo, Atlas, 2006."
$p17rng = $p17.Range
$p17rng.Find.Execute("@@PLACEHOLDER_NORMA_TOKEN@@", $true, $false, $false, $false, $false, $true, 1, $false, $bibText, 2) | Out-Null

Write-Output 'Done.'
